$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: wrap the "Brief" cell and grow the row to fit ---
$ws.Range("D16").WrapText = $true
$ws.Rows.Item(16).RowHeight = 29

# --- Row 17: new section header ---
$ws.Range("A17").Value = "1.2 V Controls"

# --- Row 18: new control (NET-003 / honeypots) ---
# Fill order mirrors the original authoring session (A, B, D, ... C filled in later)
$ws.Range("A18").Value = "NET-003"
$ws.Range("B18").Value = "Deployment of honeypots "
$ws.Range("D18").Value = "Deploy multiple honeypot files located within various directories, such as EHR folders, and fake DA accounts."

# --- Row 19: new control (AD - 006 / backup health check) ---
$ws.Range("A19").Value = "AD – 006"
$ws.Range("B19").Value = "Daily backup health check"

# back to row 18's Control Type column
$ws.Range("C18").Value = "Preventive Technical control"

# finish row 19
$ws.Range("C19").Value = "Administrative Detectice control"
$ws.Range("D19").Value = "Detects when backups are misconfigured or data integrity is compromised"

# --- Formatting for the new rows ---
$ws.Range("D18").WrapText = $true
$ws.Rows.Item(18).RowHeight = 29

$ws.Range("D19").WrapText = $true
$f19 = $ws.Range("A19").Font
$f19.Name = "Aptos"
$f19.Size = 12
$ws.Rows.Item(19).RowHeight = 29.5

# --- View state: zoomed out a bit, selection left on F18 ---
$excel.ActiveWindow.Zoom = 75
$ws.Range("F18").Select() | Out-Null
